# PowerPoint COM-interop script that reproduces the authored edit:
#  - fixes "Reijcke" -> "Rijcke" typo on the title slide
#  - tweaks bullet text on the "Inhoud presentatie" slide
#  - splits the old "Geldig vs. Goed rooster" slide's 5 bullet points into
#    5 dedicated slides ("1/5" .. "5/5"), each with its own scoring notes
#  - moves the old "Eerste stappen" and "Vragen en feedback" slides to the
#    end of the deck (with a small copy-edit on "Eerste stappen")

$p = $ppt.ActivePresentation
$layout2 = $p.SlideMaster.CustomLayouts.Item(2)   # "Titel en object" (title + content)

# ---------------------------------------------------------------------
# 1) Title slide: fix typo "Reijcke" -> "Rijcke"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleShape.TextFrame.TextRange.Replace("Reijcke", "Rijcke") | Out-Null

# ---------------------------------------------------------------------
# 2) "Inhoud presentatie" slide: annotate the "Geldig vs. Goed rooster" line
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$contentShape2 = $slide2.Shapes.Item(2)
$contentShape2.TextFrame.TextRange.Replace("Geldig vs. Goed rooster", "Geldig vs. Goed rooster (5 punten van aandacht)") | Out-Null

# ---------------------------------------------------------------------
# 3) Set aside copies of "Eerste stappen" (slide 5) and "Vragen en
#    feedback" (slide 6) at the end of the deck before we overwrite the
#    originals with the new "x/5" content.
# ---------------------------------------------------------------------
$dupEerste = $p.Slides.Item(5).Duplicate()
$dupEerste.Item(1).MoveTo($p.Slides.Count)

$dupVragen = $p.Slides.Item(6).Duplicate()
$dupVragen.Item(1).MoveTo($p.Slides.Count)

# Small copy-edit on the relocated "Eerste stappen" slide (now last-but-one)
$eerste = $p.Slides.Item($p.Slides.Count - 1)
$eersteBody = $eerste.Shapes.Item(2)
$eersteBody.TextFrame.TextRange.Replace("Vak centraal stellen zodat", "Vak centraal stellen, zodat") | Out-Null
$eersteBody.TextFrame.TextRange.Replace("Classes creëren  voor zaal, vak en student", "Classes creëren voor zaal, vak en student") | Out-Null
$eersteBody.TextFrame.TextRange.Replace("Classes creëren voor zaal, vak en student", "Classes creëren voor zaal, vak en student") | Out-Null

# ---------------------------------------------------------------------
# 4) Old slide 4 ("Geldig vs. Goed rooster") becomes the "1/5" slide
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$title4 = $slide4.Shapes.Item(1)
$body4 = $slide4.Shapes.Item(2)

$title4.TextFrame.TextRange.Text = "1/5 Alle roosterbare activiteiten hebben een tijdsslot en zaal (zaalslot)"

$body4.TextFrame.TextRange.Text = "1000 punten"
$body4.TextFrame.TextRange.Font.Size = 28

# ---------------------------------------------------------------------
# 5) Old slide 5 ("Eerste stappen") becomes the "2/5" slide
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(1)
$body5 = $slide5.Shapes.Item(2)

$title5.TextFrame.TextRange.Text = "2/5 Activiteiten per vak zoveel mogelijk verdeeld over de week"

$body5.TextFrame.TextRange.Text = "+20 per vak van 2-4 activiteiten`rVoor 2 activiteiten is dat ma-do of di-vr`rVoor 3 ma-wo-vr`rVoor 4 ma-di-do-vr`r-10 als ze op activiteiten " + [char]0x2013 + " 1 dagen geroosterd zijn`rEn -20 bij activiteiten -2 dagen et cetera"
$tr5 = $body5.TextFrame.TextRange
$tr5.Paragraphs(1,1).Font.Size = 28
$tr5.Paragraphs(2,1).IndentLevel = 2
$tr5.Paragraphs(2,1).Font.Size = 24
$tr5.Paragraphs(3,1).IndentLevel = 2
$tr5.Paragraphs(3,1).Font.Size = 24
$tr5.Paragraphs(4,1).IndentLevel = 2
$tr5.Paragraphs(4,1).Font.Size = 24
$tr5.Paragraphs(5,1).Font.Size = 28
$tr5.Paragraphs(6,1).IndentLevel = 2
$tr5.Paragraphs(6,1).Font.Size = 26

# ---------------------------------------------------------------------
# 6) Old slide 6 ("Vragen en feedback") becomes the "3/5" slide.
#    It used the title-only "Titeldia" layout; swap it for the
#    title+content layout before filling it in.
# ---------------------------------------------------------------------
$p.Slides.Item(6).Delete()
$p.Slides.AddSlide(6, $layout2) | Out-Null
$slide6 = $p.Slides.Item(6)
$slide6.Shapes.Item(1).Name = "Titel 1"
$slide6.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"

$title6 = $slide6.Shapes.Item(1)
$body6 = $slide6.Shapes.Item(2)

$title6.TextFrame.TextRange.Text = "3/5 Studenten passen in zaal"

$body6.TextFrame.TextRange.Text = "-1 per ingeschreven student die er niet inpast"
$body6.TextFrame.TextRange.Font.Size = 26

# ---------------------------------------------------------------------
# 7) Two brand new slides: "4/5" and "5/5"
# ---------------------------------------------------------------------
$p.Slides.AddSlide(7, $layout2) | Out-Null
$slide7 = $p.Slides.Item(7)
$slide7.Shapes.Item(1).Name = "Titel 1"
$slide7.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"

$title7 = $slide7.Shapes.Item(1)
$body7 = $slide7.Shapes.Item(2)
$title7.TextFrame.TextRange.Text = "4/5 Zo min mogelijk roosterconflicten per student"
$body7.TextFrame.TextRange.Text = "-1 per conflict`rConflict: meer dan " + [char]0x00E9 + [char]0x00E9 + "n activiteit in een tijdsslot"
$body7.TextFrame.TextRange.Font.Size = 26

$p.Slides.AddSlide(8, $layout2) | Out-Null
$slide8 = $p.Slides.Item(8)
$slide8.Shapes.Item(1).Name = "Titel 1"
$slide8.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"

$title8 = $slide8.Shapes.Item(1)
$body8 = $slide8.Shapes.Item(2)
$title8.TextFrame.TextRange.Text = "5/5 Werkgroepen en practica opgedeeld in zo min mogelijk groepen"
$body8.TextFrame.TextRange.Text = "Open issue"
$body8.TextFrame.TextRange.Font.Size = 26

Write-Output "Done. Final slide count: $($p.Slides.Count)"
